$wb = $excel.ActiveWorkbook

# --- Overview sheet: update status text for the b5b3f6c3 file (row 7) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E7").Value = "Handback transform failed"
$overview.Range("F7").Value = "Handback transform failed"

# --- zh-cn sheet: update status + error detail for the b5b3f6c3 file (row 7) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C7").Value = "Handback transform failed"
$zhcn.Range("P7").Value = "Handback file name: buopiqll.ibm is different with handoff file name: b5b3f6c3-da71-41d3-b427-31559db73feb.35e2f4491847f359dd16b2ad256ad3030906fa6d.zh-cn."
$zhcn.Columns.Item(16).ColumnWidth = 39.17

# --- de-de sheet: update status + error detail for the b5b3f6c3 file (row 7) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C7").Value = "Handback transform failed"
$dede.Range("P7").Value = "Handback file name: buopiqll.ibm is different with handoff file name: b5b3f6c3-da71-41d3-b427-31559db73feb.35e2f4491847f359dd16b2ad256ad3030906fa6d.de-de."
$dede.Columns.Item(16).ColumnWidth = 39.17
